$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 values (M3:W3) - replace shared formulas with literal values
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0.5
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0

# Update the selection on the active sheet to A1:X4
$ws.Range("A1:X4").Select()
